# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook tracks a changed-date value in column C for each record;
# this update advances that date by one day (46081 -> 46082, i.e.
# 2026-02-28 -> 2026-03-01) for every data row (rows 2 through 512).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 512 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 46082
